# Strip footnote markers like " [1]" and collapse embedded newlines into a
# single space across every text cell in every worksheet of the workbook.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value2

            if ($val -ne $null -and $val.GetType().Name -eq "String") {
                if ($val -match "\[\d+\]" -or $val -match "`n") {
                    $newVal = $val -replace "\[\d+\]", ""
                    $newVal = $newVal -replace "`n", " "
                    $cell.Value = $newVal
                }
            }
        }
    }
}
